$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "W01_Base_World/Simulation"
$ws.Range("C1").Value = "W01_Base_World/Simulation"

$ws.Range("C2").Select()
